$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.708.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.070.13"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.38%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.120"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.695.06"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.078.67"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.45%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.47%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.33"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.63%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0821"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.26"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "439.14"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0365"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.113"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.836.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.34"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.52"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.38%  "
